# Auto-generated from the cryptos.xlsx price/volume refresh diff.
# Updates the "Price" (D) and "Volume(1h)" (E) columns for the crypto rows
# with refreshed values, preserving each cell as plain text exactly as
# it was stored before (avoiding Excel auto-converting number-looking
# strings like "228.44" into numeric values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.534.51'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '1.810.58'
$ws.Range("E3").Value = '  +0.66%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '228.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = '  +4.00%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '34.97'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.34%  '
$ws.Range("E9").Value = '  +2.49%  '
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.0697'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.42%  '
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.0958'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.85%  '
$ws.Range("D12").Value = '2.072.40'
$ws.Range("E12").Value = '  +0.62%  '
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '11.23'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("D14").Value = '1.798.15'
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '0.651'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.80%  '
$ws.Range("E16").Value = '  +4.87%  '
$ws.Range("D17").Value = '34.516.14'
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '69.30'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.61%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '246.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("E20").Value = '  -0.50%  '
$ws.Range("E21").Value = '  +0.60%  '
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("E23").Value = '  +0.39%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '172.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.36%  '
$ws.Range("E25").Value = '  +2.19%  '
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '8.12'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +10.84%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '16.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.04%  '
$ws.Range("E28").Value = '  +2.51%  '
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("E31").Value = '  +2.19%  '
$ws.Range("E32").Value = '  +1.78%  '
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("E34").Value = '  +0.69%  '
$ws.Range("D35").Value = '1.396.49'
$ws.Range("E35").Value = '  -2.46%  '
$ws.Range("E36").Value = '  +0.69%  '
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '2.48'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.24%  '
$ws.Range("E38").Value = '  +0.68%  '
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '83.81'
$ws.Range("D40").Style = "Normal"
$ws.Range("E41").Value = '  +1.74%  '
$ws.Range("E42").Value = '  +2.72%  '
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("E44").Value = '  +5.43%  '
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '13.29'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.22%  '
$ws.Range("E46").Value = '  -2.33%  '
$ws.Range("E47").Value = '  -1.43%  '
$ws.Range("D48").Value = '1.971.95'
$ws.Range("E48").Value = '  +0.57%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '105.31'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("E50").Value = '  +1.83%  '
$ws.Range("E51").Value = '  +0.03%  '
